$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2000, "walkingToRunning", -2.787195682525635, 8.478515625, 9.063291549682615, -2.838432111238169, -4.14933462754678, 1.666932080921356),
    @(2100, "walkingToRunning", -10.50285053253174, -24.10493659973145, -11.3930606842041, -4.300542519280826, 2.241362005864094, 4.135142751430305),
    @(2200, "walkingToRunning", -2.336788654327393, 8.028789520263672, -11.69052314758301, -1.90287771350459, 13.77420247228532, -1.556544103120493),
    @(2300, "walkingToRunning", -33.8237419128418, -34.27063751220703, -13.57823276519775, -0.1651651027956011, 2.663977503776609, -2.861205784898042),
    @(2400, "walkingToRunning", -2.554568767547607, -39.26761627197266, 14.68596649169922, 8.837720880383007, -7.873765813676991, -5.242947754107039),
    @(2500, "walkingToRunning", -6.776285171508789, 9.06696891784668, -3.473462343215942, 0.8822141462250848, -0.9568282867732949, -2.844946472268377),
    @(2600, "walkingToRunning", -4.451385021209717, -37.8787841796875, 33.21585845947266, -2.166194288354044, 0.3768198490145664, 1.521733685543742),
    @(2700, "walkingToRunning", -25.62279510498047, 7.356860160827637, -15.05502891540527, -9.005332197798026, -16.56740246321032, -2.677433475067815),
    @(2800, "walkingToRunning", -8.065028190612793, -2.866869926452637, -9.438656806945801, 0.6589992062041609, -9.836473477514197, 4.410142396625722),
    @(2900, "walkingToRunning", 26.9844970703125, -14.56967449188232, 10.2812442779541, 1.44424632505363, -3.936239650375036, 3.583756440172107)
)

$startRow = 22
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}
